$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (FAPs -> FAPs) ---
$ws.Range("I2").Value = 0.7411780816145954
$ws.Range("J2").Value = 0.7411780816145954
$ws.Range("O2").Value = 0.1208074398611723
$ws.Range("P2").Value = 0.1208074398611724
$ws.Range("S2").Value = 0.08953982652107433
$ws.Range("T2").Value = 0.08953982652107434

# --- Update existing row 3 (FAPs -> MuSCs) ---
$ws.Range("I3").Value = 0.7411780816145954
$ws.Range("J3").Value = 0.7411780816145954
$ws.Range("M3").Value = 0.3894360000000001
$ws.Range("N3").Value = 1.168308
$ws.Range("O3").Value = 0.8791925601388276
$ws.Range("P3").Value = 0.8791925601388276
$ws.Range("Q3").Value = 0.018267923512
$ws.Range("R3").Value = 0.164411311608
$ws.Range("S3").Value = 0.651638255093521
$ws.Range("T3").Value = 0.651638255093521

# --- New row 4 (MuSCs -> FAPs) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Slurp1"
$ws.Range("C4").Value = "Chrna7"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01638066666666666
$ws.Range("H4").Value = 0.049142
$ws.Range("I4").Value = 0.2588219183854046
$ws.Range("J4").Value = 0.2588219183854046
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05351133333333333
$ws.Range("N4").Value = 0.160534
$ws.Range("O4").Value = 0.1208074398611723
$ws.Range("P4").Value = 0.1208074398611724
$ws.Range("Q4").Value = 0.0008765513142222222
$ws.Range("R4").Value = 0.007888961828
$ws.Range("S4").Value = 0.03126761334009803
$ws.Range("T4").Value = 0.03126761334009803

# --- New row 5 (MuSCs -> MuSCs) ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Slurp1"
$ws.Range("C5").Value = "Chrna7"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01638066666666666
$ws.Range("H5").Value = 0.049142
$ws.Range("I5").Value = 0.2588219183854046
$ws.Range("J5").Value = 0.2588219183854046
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3894360000000001
$ws.Range("N5").Value = 1.168308
$ws.Range("O5").Value = 0.8791925601388276
$ws.Range("P5").Value = 0.8791925601388276
$ws.Range("Q5").Value = 0.006379221304000001
$ws.Range("R5").Value = 0.057412991736
$ws.Range("S5").Value = 0.2275543050453066
$ws.Range("T5").Value = 0.2275543050453066
